{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the \"7.84.1\" row so the new \"7.84.2\" row can be inserted directly\n// after it (i.e. before the table's trailing empty row), regardless of\n// exactly how many rows currently exist.\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].load(\"values\");\n}\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < rows.items.length; i++) {\n  const v = rows.items[i].values;\n  if (v && v[0] && v[0][0] === \"7.84.1\") {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the '7.84.1' row to insert after.\");\n}\n\n// Insert a new (blank) row right after the \"7.84.1\" row, matching the\n// existing table formatting.\nrows.items[targetIndex].insertRows(\"After\", 1);\nawait context.sync();\n\n// Re-query the row collection to get a handle on the freshly inserted row\n// (index targetIndex + 1).\nrows.load(\"items\");\nawait context.sync();\n\nconst newRow = rows.items[targetIndex + 1];\n\nconst newCells = newRow.cells;\nnewCells.load(\"items\");\nawait context.sync();\n\nconst revisionXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  \"<w:p><w:r><w:t>7.84.2</w:t></w:r></w:p>\" +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst notesXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  \"<w:p><w:r><w:t>Common ITs:</w:t></w:r></w:p>\" +\n  \"<w:p><w:pPr><w:ind w:left=\\\"720\\\"/></w:pPr>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\">#357 \\u2013 Add </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>hwReadAndCheckByte</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\">() and </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>hwReadAndCheckWord</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"<w:r><w:t>() to hw_intf.c</w:t></w:r></w:p>\" +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nnewCells.items[0].getRange().insertOoxml(revisionXml, Word.InsertLocation.replace);\nnewCells.items[1].getRange().insertOoxml(notesXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Find the \"7.84.1\" row so the new \"7.84.2\" row can be inserted directly\n# after it (i.e. before the table's trailing empty row), regardless of\n# exactly how many rows currently exist.\n$targetRowIndex = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n  $cellText = $t.Cell($i, 1).Range.Text\n  $cellText = $cellText.TrimEnd([char]13, [char]7)\n  if ($cellText -eq \"7.84.1\") {\n    $targetRowIndex = $i\n    break\n  }\n}\nif ($targetRowIndex -eq -1) {\n  throw \"Could not find the '7.84.1' row to insert after.\"\n}\n\n# Insert a new (blank) row right after the \"7.84.1\" row, matching the\n# table's existing formatting (Rows.Add clones the row immediately above\n# the row passed in, i.e. the row it is inserted before).\n$beforeRow = $t.Rows.Item($targetRowIndex + 1)\n$newRow = $t.Rows.Add($beforeRow)\n$newRowIndex = $targetRowIndex + 1\n\n# --- Revision cell (\"7.84.2\") ---\n$revisionXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>7.84.2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n$t.Cell($newRowIndex, 1).Range.InsertXML($revisionXml) | Out-Null\n\n# --- Notes cell (\"Common ITs: / #357 ...\") ---\n$notesXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Common ITs:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">#357 \u2013 Add </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hwReadAndCheckByte</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">() and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hwReadAndCheckWord</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>() to hw_intf.c</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n$t.Cell($newRowIndex, 2).Range.InsertXML($notesXml) | Out-Null\n"}
